$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.113.50'
$ws.Range('E2').Value = '  +6.42%  '
$ws.Range('D3').Value = '1.895.18'
$ws.Range('E3').Value = '  +6.03%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.9992'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '250.06'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.40%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.9992'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.01%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5010'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +1.87%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '45.96'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +9.06%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.2880'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +6.96%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.06579'
$cell.Style = 'Normal'
$ws.Range('D11').Value = '1.885.04'
$ws.Range('E11').Value = '  +5.45%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '17.31'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +4.54%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.07252'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('E14').Value = '  +6.83%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '85.12'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +6.24%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '4.842'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +3.84%  '
$ws.Range('D17').Value = '30.129.14'
$ws.Range('E17').Value = '  +6.66%  '
$ws.Range('E18').Value = '  -0.03%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '12.95'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +7.29%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.000007549'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +4.07%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '0.9989'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = '2.128.62'
$ws.Range('E22').Value = '  +5.66%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '4.791'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +5.06%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '5.567'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +5.96%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '9.069'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +3.32%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '145.71'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.61%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '136.63'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +24.38%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '16.84'
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.957'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +5.60%  '
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E31').Value = '  +1.17%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.08690'
$cell.Style = 'Normal'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '3.954'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +4.77%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.05040'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +3.06%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.143'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('E36').Value = '  +5.99%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.687'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('E38').Value = '  +10.76%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.779'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +6.43%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.9646'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +1.84%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.01644'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +5.90%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '6.072'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +2.65%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '105.40'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +5.54%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.9990'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.4234'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +5.89%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '7.482'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +3.98%  '
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('E48').Value = '  +3.94%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '32.63'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +6.02%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '8.321'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +3.60%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.3741'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +6.94%  '
